$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Negative words table (A3:H6) ---
$neg = New-Object 'object[,]' 4,8
$neg[0,0] = "crude"
$neg[0,1] = 0.8823529411764706
$neg[0,2] = 30
$neg[0,3] = 30
$neg[0,4] = 0
$neg[0,5] = 1
$neg[0,6] = $false
$neg[0,7] = 4
$neg[1,0] = "crisis"
$neg[1,1] = 0.5856164383561644
$neg[1,2] = 171
$neg[1,3] = 171
$neg[1,4] = 0
$neg[1,5] = 1
$neg[1,6] = $false
$neg[1,7] = 121
$neg[2,0] = "sc"
$neg[2,1] = 0.1904761904761905
$neg[2,2] = 36
$neg[2,3] = 36
$neg[2,4] = 0
$neg[2,5] = 1
$neg[2,6] = $false
$neg[2,7] = 153
$neg[3,0] = "panic"
$neg[3,1] = 0.1724806201550388
$neg[3,2] = 89
$neg[3,3] = 89
$neg[3,4] = 0
$neg[3,5] = 1
$neg[3,6] = $false
$neg[3,7] = 427
$ws.Range("A3:H6").Value = $neg

# Remove the old 5th negative-word row (was "low"); shift remaining cells up so the row disappears cleanly
$ws.Range("A7:H7").Delete(-4162)

# --- Positive words table (J3:Q27) ---
$pos = New-Object 'object[,]' 25,8
$pos[0,0] = "best"
$pos[0,1] = 0.9491525423728814
$pos[0,2] = 56
$pos[0,3] = 56
$pos[0,4] = 1
$pos[0,5] = 0
$pos[0,6] = $false
$pos[0,7] = 3
$pos[1,0] = "interesting"
$pos[1,1] = 0.9393939393939394
$pos[1,2] = 31
$pos[1,3] = 31
$pos[1,4] = 1
$pos[1,5] = 0
$pos[1,6] = $false
$pos[1,7] = 2
$pos[2,0] = "love"
$pos[2,1] = 0.9130434782608695
$pos[2,2] = 42
$pos[2,3] = 42
$pos[2,4] = 1
$pos[2,5] = 0
$pos[2,6] = $false
$pos[2,7] = 4
$pos[3,0] = "great"
$pos[3,1] = 0.8660714285714286
$pos[3,2] = 97
$pos[3,3] = 97
$pos[3,4] = 1
$pos[3,5] = 0
$pos[3,6] = $false
$pos[3,7] = 15
$pos[4,0] = "positive"
$pos[4,1] = 0.8275862068965517
$pos[4,2] = 48
$pos[4,3] = 48
$pos[4,4] = 1
$pos[4,5] = 0
$pos[4,6] = $false
$pos[4,7] = 10
$pos[5,0] = "thanks"
$pos[5,1] = 0.8170731707317073
$pos[5,2] = 67
$pos[5,3] = 67
$pos[5,4] = 1
$pos[5,5] = 0
$pos[5,6] = $false
$pos[5,7] = 15
$pos[6,0] = "special"
$pos[6,1] = 0.8055555555555556
$pos[6,2] = 29
$pos[6,3] = 29
$pos[6,4] = 1
$pos[6,5] = 0
$pos[6,6] = $false
$pos[6,7] = 7
$pos[7,0] = "thank"
$pos[7,1] = 0.78125
$pos[7,2] = 100
$pos[7,3] = 100
$pos[7,4] = 1
$pos[7,5] = 0
$pos[7,6] = $false
$pos[7,7] = 28
$pos[8,0] = "safety"
$pos[8,1] = 0.7254901960784313
$pos[8,2] = 37
$pos[8,3] = 37
$pos[8,4] = 1
$pos[8,5] = 0
$pos[8,6] = $false
$pos[8,7] = 14
$pos[9,0] = "free"
$pos[9,1] = 0.7166666666666667
$pos[9,2] = 86
$pos[9,3] = 86
$pos[9,4] = 1
$pos[9,5] = 0
$pos[9,6] = $false
$pos[9,7] = 34
$pos[10,0] = "confidence"
$pos[10,1] = 0.6944444444444444
$pos[10,2] = 25
$pos[10,3] = 25
$pos[10,4] = 1
$pos[10,5] = 0
$pos[10,6] = $false
$pos[10,7] = 11
$pos[11,0] = "safe"
$pos[11,1] = 0.6901408450704225
$pos[11,2] = 98
$pos[11,3] = 98
$pos[11,4] = 1
$pos[11,5] = 0
$pos[11,6] = $false
$pos[11,7] = 44
$pos[12,0] = "support"
$pos[12,1] = 0.6886792452830188
$pos[12,2] = 73
$pos[12,3] = 73
$pos[12,4] = 1
$pos[12,5] = 0
$pos[12,6] = $false
$pos[12,7] = 33
$pos[13,0] = "good"
$pos[13,1] = 0.6625
$pos[13,2] = 106
$pos[13,3] = 106
$pos[13,4] = 1
$pos[13,5] = 0
$pos[13,6] = $false
$pos[13,7] = 54
$pos[14,0] = "relief"
$pos[14,1] = 0.62
$pos[14,2] = 31
$pos[14,3] = 31
$pos[14,4] = 1
$pos[14,5] = 0
$pos[14,6] = $false
$pos[14,7] = 19
$pos[15,0] = "well"
$pos[15,1] = 0.5957446808510638
$pos[15,2] = 56
$pos[15,3] = 56
$pos[15,4] = 1
$pos[15,5] = 0
$pos[15,6] = $false
$pos[15,7] = 38
$pos[16,0] = "heroes"
$pos[16,1] = 0.5319148936170213
$pos[16,2] = 25
$pos[16,3] = 25
$pos[16,4] = 1
$pos[16,5] = 0
$pos[16,6] = $false
$pos[16,7] = 22
$pos[17,0] = "better"
$pos[17,1] = 0.5238095238095238
$pos[17,2] = 33
$pos[17,3] = 33
$pos[17,4] = 1
$pos[17,5] = 0
$pos[17,6] = $false
$pos[17,7] = 30
$pos[18,0] = "hand"
$pos[18,1] = 0.5143603133159269
$pos[18,2] = 197
$pos[18,3] = 197
$pos[18,4] = 1
$pos[18,5] = 0
$pos[18,6] = $false
$pos[18,7] = 186
$pos[19,0] = "care"
$pos[19,1] = 0.4831460674157304
$pos[19,2] = 43
$pos[19,3] = 43
$pos[19,4] = 1
$pos[19,5] = 0
$pos[19,6] = $false
$pos[19,7] = 46
$pos[20,0] = "like"
$pos[20,1] = 0.4470588235294118
$pos[20,2] = 152
$pos[20,3] = 152
$pos[20,4] = 1
$pos[20,5] = 0
$pos[20,6] = $false
$pos[20,7] = 188
$pos[21,0] = "help"
$pos[21,1] = 0.4135593220338983
$pos[21,2] = 122
$pos[21,3] = 122
$pos[21,4] = 1
$pos[21,5] = 0
$pos[21,6] = $false
$pos[21,7] = 173
$pos[22,0] = "increase"
$pos[22,1] = 0.358974358974359
$pos[22,2] = 28
$pos[22,3] = 28
$pos[22,4] = 1
$pos[22,5] = 0
$pos[22,6] = $false
$pos[22,7] = 50
$pos[23,0] = "protect"
$pos[23,1] = 0.3561643835616438
$pos[23,2] = 26
$pos[23,3] = 26
$pos[23,4] = 1
$pos[23,5] = 0
$pos[23,6] = $false
$pos[23,7] = 47
$pos[24,0] = "please"
$pos[24,1] = 0.3389121338912134
$pos[24,2] = 81
$pos[24,3] = 81
$pos[24,4] = 1
$pos[24,5] = 0
$pos[24,6] = $false
$pos[24,7] = 158
$ws.Range("J3:Q27").Value = $pos
